$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new monthly values (Julho..Novembro 2023) to row 2
$ws.Range("H2").Value = 1.0720000000000001
$ws.Range("I2").Value = 1.1375
$ws.Range("J2").Value = 0.97289999999999999
$ws.Range("K2").Value = 0.99760000000000004
$ws.Range("L2").Value = 0.91600000000000004

# Move the "Fonte" hyperlink from O3 up to O2, copying its text + style.
$hlStyle = $ws.Range("O3").Style

$ws.Range("O3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("O2"), "https://www.valor.srv.br/indices/cdi.php")

$ws.Range("O2").Value = "https://www.valor.srv.br/indices/cdi.php"
$ws.Range("O2").Style = $hlStyle

# Clear out the now-empty O3 cell
$ws.Range("O3").Clear()
